# edit.ps1 - reproduce the commit:
#   1. Slide 5's table switches from the custom "Table_0" style to the
#      built-in "No Style, Table Grid" style.
#   2. The theme applied via the Slide Master (ppt/theme/theme1.xml,
#      previously "Integral" / "Red Violet") is recoloured to match the
#      stock "Office Theme" colour scheme (the colours that used to live,
#      unused, in the Notes Master's theme part).

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{B14F1FB9-2B58-4C1F-A4DC-D234AA0E998F}")

# --- 2. Theme colour scheme -------------------------------------------------
function HexToComRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# msoThemeColorSchemeIndex slot -> target "Office Theme" colour.
# (dark1/light1 are already 000000/FFFFFF in both themes, so only the
# remaining ten slots actually change.)
$officeColors = @{
    3  = "44546A"  # Dark 2
    4  = "E7E6E6"  # Light 2
    5  = "5B9BD5"  # Accent 1
    6  = "ED7D31"  # Accent 2
    7  = "A5A5A5"  # Accent 3
    8  = "FFC000"  # Accent 4
    9  = "4472C4"  # Accent 5
    10 = "70AD47"  # Accent 6
    11 = "0563C1"  # Hyperlink
    12 = "954F72"  # Followed Hyperlink
}

$master = $p.Slides.Item(1).Master
$colorScheme = $master.Theme.ThemeColorScheme

foreach ($index in $officeColors.Keys) {
    $colorScheme.Item($index).RGB = HexToComRgb $officeColors[$index]
}
